$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3  = @(301, 6, 45, 30, 60, 45)
    4  = @(1203, 3, 15, 15, 15, 15)
    5  = @(101, 9, 30, 15, 60, 15)
    6  = @(801, 3, 67, 65, 52, 45)
    7  = @(1202, 2, 10, 10, 10, 10)
    10 = @(601, 9, 60, 67, 60, 42)
    11 = @(902, 1, 0, 0, 0, 0)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(701, 3, 90, 45, 97, 15)
    14 = @(201, 9, 30, 15, 45, 30)
    15 = @(1201, 2, 10, 10, 10, 10)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(802, 0, 4, 5, 4, 0)
    18 = @(3, 0, 3, 3, 3, 3)
    19 = @(1101, 0, 15, 30, 30, 0)
    20 = @(1, 0, 2, 2, 2, 2)
    21 = @(2, 0, 2, 2, 2, 2)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
